$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "34.482.69"
$ws.Cells.Item(2, 5).Value = "  +0.38%  "
$ws.Cells.Item(3, 4).Value = "1.813.62"
$ws.Cells.Item(3, 5).Value = "  +0.65%  "
$ws.Cells.Item(4, 5).Value = "  -0.21%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "225.94"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.63%  "
$ws.Cells.Item(7, 5).Value = "  -0.20%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "38.42"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +7.04%  "
$ws.Cells.Item(9, 5).Value = "  -3.62%  "
$ws.Cells.Item(10, 5).Value = "  -2.34%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0974"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +1.08%  "
$ws.Cells.Item(12, 4).Value = "2.074.51"
$ws.Cells.Item(12, 5).Value = "  +0.56%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "11.24"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -2.30%  "
$ws.Cells.Item(14, 4).Value = "1.829.57"
$ws.Cells.Item(14, 5).Value = "  +1.57%  "
$ws.Cells.Item(15, 5).Value = "  -1.62%  "
$ws.Cells.Item(16, 4).Value = "34.462.17"
$ws.Cells.Item(16, 5).Value = "  +0.29%  "
$ws.Cells.Item(17, 5).Value = "  -1.53%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "68.37"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -0.95%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "243.19"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -1.00%  "
$ws.Cells.Item(20, 4).Value = "0.0₃0773"
$ws.Cells.Item(20, 5).Value = "  -2.61%  "
$ws.Cells.Item(21, 5).Value = "  -1.97%  "
$ws.Cells.Item(22, 5).Value = "  -0.17%  "
$ws.Cells.Item(23, 5).Value = "  -0.97%  "
$ws.Cells.Item(24, 5).Value = "  +3.52%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "170.18"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -0.37%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "7.83"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -0.36%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "17.60"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +4.38%  "
$ws.Cells.Item(28, 5).Value = "  +1.89%  "
$ws.Cells.Item(29, 5).Value = "  -0.18%  "
$ws.Cells.Item(30, 2).Value = "Filecoin"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "3.80"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -0.98%  "
$ws.Cells.Item(31, 2).Value = "PancakeSwap"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.23"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -0.90%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.0519"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -2.30%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "3.86"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -4.44%  "
$ws.Cells.Item(35, 4).Value = "1.362.18"
$ws.Cells.Item(35, 5).Value = "  -2.21%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.645"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -4.03%  "
$ws.Cells.Item(37, 5).Value = "  -0.32%  "
$ws.Cells.Item(38, 5).Value = "  -4.29%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.45"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +1.38%  "
$ws.Cells.Item(41, 5).Value = "  -0.94%  "
$ws.Cells.Item(42, 2).Value = "Aave"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "81.88"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -0.20%  "
$ws.Cells.Item(43, 2).Value = "ARBITRUM"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.952"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -0.91%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "2.81"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -0.58%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "13.78"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +1.69%  "
$ws.Cells.Item(46, 5).Value = "  +1.18%  "
$ws.Cells.Item(47, 4).Value = "1.975.60"
$ws.Cells.Item(47, 5).Value = "  +0.58%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "5.77"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -4.50%  "
$ws.Cells.Item(49, 5).Value = "  -0.17%  "
$ws.Cells.Item(50, 5).Value = "  -2.49%  "
$ws.Cells.Item(51, 5).Value = "  -4.62%  "
